$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.975.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.782.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.12%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0709"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.039.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.780.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.947.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0522"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.399.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "79.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0491"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.937.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("E51").Value = "  -1.40%  "
